# Regenerate the "K" (strikeouts) column (column G) in the save_data sheet
# for kopech_michael.xlsx. This replaces the old "Strike#" derived values
# with the freshly-computed K values (std/mean regenerated upstream, s_vals
# written here), rows 2-53 correspond to data rows 0-51 in column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values, in row order starting at row 2 through row 53.
$newK = @(0,5,6,3,0,4,1,2,3,5,2,3,2,2,1,1,1,3,2,1,0,1,1,0,3,0,4,1,3,3,0,0,0,3,2,0,5,1,4,3,10,4,2,3,5,4,1,1,2,2,2,1)

$startRow = 2
for ($i = 0; $i -lt $newK.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
